$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2331.6667
$ws.Range("J17").Value = 2997.5
$ws.Range("L17").Value = 8992.5
$ws.Range("N17").Value = -9328.5
$ws.Range("H53").Value = 251.33333
$ws.Range("I53").Value = 309.83334
$ws.Range("J53").Value = 134.33333
$ws.Range("K53").Value = 309.83334
$ws.Range("L53").Value = 134.33333
$ws.Range("M53").Value = 327.16666
$ws.Range("N53").Value = -1408.33333
$ws.Range("H58").Value = 1680
$ws.Range("I58").Value = 193.33333
$ws.Range("J58").Value = 3166.6667
$ws.Range("K58").Value = 579.99999
$ws.Range("L58").Value = 9500.000100000001
$ws.Range("M58").Value = -429.99999
$ws.Range("N58").Value = -9800.000100000001
$ws.Range("H74").Value = 3983.625
$ws.Range("I74").Value = 3983.625
$ws.Range("K74").Value = 3983.625
$ws.Range("M74").Value = -3047.625
$ws.Range("H77").Value = 3983.625
$ws.Range("I77").Value = 3983.625
$ws.Range("K77").Value = 19918.125
$ws.Range("M77").Value = -15238.125
$ws.Range("H141").Value = 11439.4
$ws.Range("I141").Value = 13948
$ws.Range("K141").Value = 41844
$ws.Range("M141").Value = -36664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3756.25
$ws.Range("I45").Value = 3003.6667
$ws.Range("K45").Value = 3003.6667
$ws.Range("M45").Value = -2626.6667
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H132").Value = 1737.125
$ws.Range("I132").Value = 1732.9333
$ws.Range("K132").Value = 5198.7999
$ws.Range("M132").Value = -2668.7999
$ws.Range("H135").Value = 36215
$ws.Range("J135").Value = 36215
$ws.Range("L135").Value = 36215
$ws.Range("N135").Value = -46355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 497.4737
$ws.Range("J7").Value = 608.2222
$ws.Range("L7").Value = 608.2222
$ws.Range("N7").Value = -834.2222
$ws.Range("H31").Value = 3599.5386
$ws.Range("I31").Value = 1500
$ws.Range("J31").Value = 3981.2727
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 3981.2727
$ws.Range("M31").Value = -1205
$ws.Range("N31").Value = -4571.2727
$ws.Range("H34").Value = 3599.5386
$ws.Range("I34").Value = 1500
$ws.Range("J34").Value = 3981.2727
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 3981.2727
$ws.Range("M34").Value = -1298
$ws.Range("N34").Value = -4385.2727
$ws.Range("H99").Value = 4150.44
$ws.Range("I99").Value = 4261.864
$ws.Range("K99").Value = 4261.864
$ws.Range("M99").Value = -2763.864
$ws.Range("H118").Value = 84999.5
$ws.Range("J118").Value = 84999.5
$ws.Range("L118").Value = 84999.5
$ws.Range("N118").Value = -88313.5
$ws.Range("H126").Value = 4150.44
$ws.Range("I126").Value = 4261.864
$ws.Range("K126").Value = 12785.592
$ws.Range("M126").Value = -10315.592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2974.5264
$ws.Range("I34").Value = 310.4
$ws.Range("J34").Value = 3926
$ws.Range("K34").Value = 931.1999999999999
$ws.Range("L34").Value = 11778
$ws.Range("M34").Value = -847.1999999999999
$ws.Range("N34").Value = -11946
$ws.Range("H94").Value = 20800
$ws.Range("J94").Value = 20800
$ws.Range("L94").Value = 62400
$ws.Range("N94").Value = -63752
$ws.Range("H121").Value = 627120.2
$ws.Range("J121").Value = 1112354.9
$ws.Range("L121").Value = 3337064.7
$ws.Range("N121").Value = -3339684.7
$ws.Range("H139").Value = 9258.875
$ws.Range("I139").Value = 9258.875
$ws.Range("K139").Value = 27776.625
$ws.Range("M139").Value = -22636.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4335
$ws.Range("I80").Value = 3005
$ws.Range("K80").Value = 3005
$ws.Range("M80").Value = -2007
$ws.Range("H83").Value = 4335
$ws.Range("I83").Value = 3005
$ws.Range("K83").Value = 15025
$ws.Range("M83").Value = -10033
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2944
$ws.Range("I40").Value = 2940.7273
$ws.Range("J40").Value = 2956
$ws.Range("K40").Value = 2940.7273
$ws.Range("L40").Value = 2956
$ws.Range("M40").Value = -2804.7273
$ws.Range("N40").Value = -3228
$ws.Range("H46").Value = 1471.3334
$ws.Range("I46").Value = 1599
$ws.Range("K46").Value = 1599
$ws.Range("M46").Value = -1411
$ws.Range("H130").Value = 82164.5
$ws.Range("J130").Value = 82164.5
$ws.Range("L130").Value = 82164.5
$ws.Range("N130").Value = -92204.5
$ws.Range("H132").Value = 3042.4375
$ws.Range("I132").Value = 3084.2144
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 9252.643199999999
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -6722.643199999999
$ws.Range("N132").Value = -13310
$ws.Range("H136").Value = 20092.777
$ws.Range("J136").Value = 21731.5
$ws.Range("L136").Value = 65194.5
$ws.Range("N136").Value = -70294.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 69000
$ws.Range("J47").Value = 69000
$ws.Range("L47").Value = 69000
$ws.Range("N47").Value = -70144
$ws.Range("H51").Value = 201111.5
$ws.Range("I51").Value = 256665.33
$ws.Range("J51").Value = 34450
$ws.Range("K51").Value = 256665.33
$ws.Range("L51").Value = 34450
$ws.Range("M51").Value = -256155.33
$ws.Range("N51").Value = -35470
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("J56").Value = 15000
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -14286
$ws.Range("N56").Value = -16428
$ws.Range("H62").Value = 4999.5
$ws.Range("J62").Value = 4999.5
$ws.Range("L62").Value = 4999.5
$ws.Range("N62").Value = -6247.5
$ws.Range("H65").Value = 4999.5
$ws.Range("J65").Value = 4999.5
$ws.Range("L65").Value = 24997.5
$ws.Range("N65").Value = -31237.5
$ws.Range("H81").Value = 3005.7856
$ws.Range("I81").Value = 1916.5454
$ws.Range("J81").Value = 6999.6665
$ws.Range("K81").Value = 3833.0908
$ws.Range("L81").Value = 13999.333
$ws.Range("M81").Value = -2772.0908
$ws.Range("N81").Value = -16121.333
$ws.Range("H84").Value = 3005.7856
$ws.Range("I84").Value = 1916.5454
$ws.Range("J84").Value = 6999.6665
$ws.Range("K84").Value = 19165.454
$ws.Range("L84").Value = 69996.66500000001
$ws.Range("M84").Value = -13861.454
$ws.Range("N84").Value = -80604.66500000001
$ws.Range("H122").Value = 2750.3845
$ws.Range("I122").Value = 2562.5
$ws.Range("J122").Value = 5005
$ws.Range("K122").Value = 7687.5
$ws.Range("L122").Value = 15015
$ws.Range("M122").Value = -5237.5
$ws.Range("N122").Value = -19915
$ws.Range("H126").Value = 2400.6843
$ws.Range("I126").Value = 2816.3845
$ws.Range("K126").Value = 8449.1535
$ws.Range("M126").Value = -5979.1535
$ws.Range("H136").Value = 15868.917
$ws.Range("I136").Value = 12652.8
$ws.Range("J136").Value = 31949.5
$ws.Range("K136").Value = 37958.39999999999
$ws.Range("L136").Value = 95848.5
$ws.Range("M136").Value = -35408.39999999999
$ws.Range("N136").Value = -100948.5
